$d = $word.ActiveDocument

# Update the header date
$d.Content.Find.Execute("2024-06-19 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-06-20 Thursday", 2)

# Update the division problems in the worksheet table, row by row (1-based
# row indices correspond to the 5 populated rows: 1, 5, 9, 13, 17).
$t = $d.Tables.Item(1)

$rowsData = @{
    1  = @("930÷7=", "370÷6=", "655÷3=", "675÷2=", "351÷4=")
    5  = @("676÷3=", "309÷6=", "312÷7=", "399÷2=", "592÷8=")
    9  = @("461÷7=", "842÷6=", "479÷4=", "340÷5=", "400÷8=")
    13 = @("522÷6=", "218÷4=", "542÷6=", "629÷8=", "858÷5=")
    17 = @("246÷5=", "185÷8=", "246÷8=", "345÷5=", "874÷8=")
}

foreach ($rowIndex in $rowsData.Keys) {
    $row = $t.Rows.Item($rowIndex)
    $values = $rowsData[$rowIndex]
    for ($col = 1; $col -le $values.Count; $col++) {
        $row.Cells.Item($col).Range.Text = $values[$col - 1]
    }
}
